# Fruta / hortaliza, semanal
# Insert one new weekly pair of price records (Pintón + Primera Pintón) for
# Plátano at the Macroferia Regional de Talca market. The new records are
# inserted right above the existing row 487, which pushes all subsequent
# rows down by two rows (487->489, ... 581->583) exactly like a normal
# Excel row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 487 (shifts 487:581 down to 489:583)
$ws.Rows("487:488").Insert()

# The row that used to be 487 (template for quality "Pintón") is now at 489;
# the row that used to be 488 (template for quality "Primera Pintón") is now
# at 490. Copy their full contents (formats + all the constant columns such
# as Mercado, Región, Producto, Unidad, Origen, etc.) into the two freshly
# inserted blank rows so every column besides the ones that actually change
# is populated exactly like the rest of the table.
$ws.Range("A489:T489").Copy()
$ws.Range("A487:T487").PasteSpecial()

$ws.Range("A490:T490").Copy()
$ws.Range("A488:T488").PasteSpecial()

# Now overwrite the cells that hold the new week's data.
# Row 487: Calidad "Pintón"
$ws.Range("D487").Value = 44641
$ws.Range("M487").Value = 1050
$ws.Range("N487:P487").Value = 18500
$ws.Range("S487").Value = 925

# Row 488: Calidad "Primera Pintón"
$ws.Range("D488").Value = 44641
$ws.Range("M488").Value = 480
$ws.Range("N488:P488").Value = 19500
$ws.Range("S488").Value = 975
